$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scaling Parameter")

$ws.Range("C2").Value = 0.2883263863596509
$ws.Range("C3").Value = 0.7214736462234707
$ws.Range("C4").Value = 0.5390408112743797
$ws.Range("C5").Value = 0.2551466653460389
$ws.Range("C6").Value = 0.1629806062963671
$ws.Range("C7").Value = 0.00218285291877274
$ws.Range("C8").Value = 0.1549130199603967
$ws.Range("C9").Value = 0.07328113654965122
$ws.Range("C10").Value = 0.2925471910957358
$ws.Range("C11").Value = 0.02482950365147664
$ws.Range("C12").Value = 0.07562380854834035
$ws.Range("C13").Value = 0.5994515917325207
$ws.Range("C14").Value = 0.049038861786753
$ws.Range("C15").Value = 0.06552336369281705
$ws.Range("C16").Value = 0.01503284782456222
$ws.Range("C17").Value = 0.00419105407513218
$ws.Range("C18").Value = 0.03458551375137886
$ws.Range("C19").Value = 0.00211206463118908
$ws.Range("C20").Value = 0.00003034357094102027
$ws.Range("C21").Value = 0.0177020890977593
$ws.Range("C22").Value = 0.9000002546482861
$ws.Range("C23").Value = 0.2298898440478568
$ws.Range("C24").Value = 0.2583499585106195
$ws.Range("C25").Value = 0.1747826671811203
$ws.Range("C26").Value = 0.2226978885588895
$ws.Range("C27").Value = 0.1827788481342695
$ws.Range("C28").Value = 0.1779914587157404
$ws.Range("C29").Value = 0.04889878998626204
$ws.Range("C30").Value = 0.1442166415905932
$ws.Range("C31").Value = 0.2434215849764926
$ws.Range("C32").Value = 0.2190527543621555
$ws.Range("C33").Value = 0.2583500250428513
$ws.Range("C34").Value = 0.258350551527436
$ws.Range("C35").Value = 0.2417665505285693
$ws.Range("C36").Value = 0.2231959418462669
